# Auto-generated edit script applying price/profit data updates
# as per the scheduled runner data refresh (see diff).
$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 35.714287
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 10
$ws.Range("K5").Value = 70
$ws.Range("L5").Value = 10
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = -240

# Row 33
$ws.Range("H33").Value = 213.9
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

# Row 43
$ws.Range("H43").Value = 33335378
$ws.Range("I43").Value = 125003120
$ws.Range("J43").Value = 1654.5454
$ws.Range("K43").Value = 125003120
$ws.Range("L43").Value = 1654.5454
$ws.Range("M43").Value = -125003051
$ws.Range("N43").Value = -1792.5454

# Row 53
$ws.Range("H53").Value = 16129532
$ws.Range("I53").Value = 35714436
$ws.Range("J53").Value = 788.82355
$ws.Range("K53").Value = 35714436
$ws.Range("L53").Value = 788.82355
$ws.Range("M53").Value = -35713799
$ws.Range("N53").Value = -2062.82355

# Row 55
$ws.Range("H55").Value = 1202403
$ws.Range("I55").Value = 826.6923
$ws.Range("J55").Value = 2403979.2
$ws.Range("K55").Value = 826.6923
$ws.Range("L55").Value = 2403979.2
$ws.Range("M55").Value = -612.6923
$ws.Range("N55").Value = -2404407.2

# Row 69
$ws.Range("H69").Value = 1750
$ws.Range("I69").Value = 1750
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 5250
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -4376
$ws.Range("N69").ClearContents()

# Row 72
$ws.Range("H72").Value = 1750
$ws.Range("I72").Value = 1750
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 15750
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -11382
$ws.Range("N72").ClearContents()

# Row 80
$ws.Range("H80").Value = 443
$ws.Range("I80").Value = 374.75
$ws.Range("J80").Value = 534
$ws.Range("K80").Value = 1124.25
$ws.Range("L80").Value = 1602
$ws.Range("M80").Value = -126.25
$ws.Range("N80").Value = -3598

# Row 83
$ws.Range("H83").Value = 443
$ws.Range("I83").Value = 374.75
$ws.Range("J83").Value = 534
$ws.Range("K83").Value = 3372.75
$ws.Range("L83").Value = 4806
$ws.Range("M83").Value = 1619.25
$ws.Range("N83").Value = -14790

# Row 132
$ws.Range("H132").Value = 5095.5576
$ws.Range("I132").Value = 4169.418
$ws.Range("K132").Value = 12508.254
$ws.Range("M132").Value = -9978.253999999999

# Row 137
$ws.Range("H137").Value = 1080.375
$ws.Range("I137").Value = 848.1707
$ws.Range("J137").Value = 2440.4285
$ws.Range("K137").Value = 2544.5121
$ws.Range("L137").Value = 7321.2855
$ws.Range("M137").Value = 5.487900000000081
$ws.Range("N137").Value = -12421.2855

# Row 138
$ws.Range("H138").Value = 2038.4067
$ws.Range("I138").Value = 904.25
$ws.Range("J138").Value = 3382.5925
$ws.Range("K138").Value = 2712.75
$ws.Range("L138").Value = 10147.7775
$ws.Range("M138").Value = 2427.25
$ws.Range("N138").Value = -20427.7775


# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1242.6
$ws.Range("I2").Value = 1450
$ws.Range("J2").Value = 1104.3334
$ws.Range("K2").Value = 1450
$ws.Range("L2").Value = 1104.3334
$ws.Range("M2").Value = -1337
$ws.Range("N2").Value = -1330.3334

# Row 32
$ws.Range("H32").Value = 16399062
$ws.Range("I32").Value = 4863.4365
$ws.Range("K32").Value = 4863.4365
$ws.Range("M32").Value = -4576.4365

# Row 102
$ws.Range("H102").Value = 1255.5
$ws.Range("I102").Value = 500
$ws.Range("K102").Value = 500
$ws.Range("M102").Value = 1122

# Row 116
$ws.Range("H116").Value = 1242.6
$ws.Range("I116").Value = 1450
$ws.Range("J116").Value = 1104.3334
$ws.Range("K116").Value = 1450
$ws.Range("L116").Value = 1104.3334
$ws.Range("M116").Value = 844
$ws.Range("N116").Value = -5692.3334

# Row 122
$ws.Range("H122").Value = 1522.5
$ws.Range("I122").Value = 1284.375
$ws.Range("J122").Value = 1998.75
$ws.Range("K122").Value = 3853.125
$ws.Range("L122").Value = 5996.25
$ws.Range("M122").Value = -1403.125
$ws.Range("N122").Value = -10896.25

# Row 132
$ws.Range("H132").Value = 1033108.56
$ws.Range("I132").Value = 1029.7675
$ws.Range("J132").Value = 4203065
$ws.Range("K132").Value = 3089.3025
$ws.Range("L132").Value = 12609195
$ws.Range("M132").Value = -559.3024999999998
$ws.Range("N132").Value = -12614255


# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1242.6
$ws.Range("I3").Value = 1450
$ws.Range("J3").Value = 1104.3334
$ws.Range("K3").Value = 1450
$ws.Range("L3").Value = 1104.3334
$ws.Range("M3").Value = -1336
$ws.Range("N3").Value = -1332.3334

# Row 22
$ws.Range("H22").Value = 6757006.5
$ws.Range("I22").Value = 6757006.5
$ws.Range("K22").Value = 6757006.5
$ws.Range("M22").Value = -6756833.5

# Row 105
$ws.Range("H105").Value = 52633220
$ws.Range("I105").Value = 1657.7142
$ws.Range("K105").Value = 1657.7142
$ws.Range("M105").Value = 89.28580000000011


# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1559.2142
$ws.Range("I16").Value = 1102.6364
$ws.Range("J16").Value = 3233.3333
$ws.Range("K16").Value = 1102.6364
$ws.Range("L16").Value = 3233.3333
$ws.Range("M16").Value = -815.6364000000001
$ws.Range("N16").Value = -3807.3333

# Row 31
$ws.Range("H31").Value = 2348.342
$ws.Range("I31").Value = 2512.8813
$ws.Range("J31").Value = 1777.2941
$ws.Range("K31").Value = 2512.8813
$ws.Range("L31").Value = 1777.2941
$ws.Range("M31").Value = -2217.8813
$ws.Range("N31").Value = -2367.2941

# Row 34
$ws.Range("H34").Value = 2348.342
$ws.Range("I34").Value = 2512.8813
$ws.Range("J34").Value = 1777.2941
$ws.Range("K34").Value = 2512.8813
$ws.Range("L34").Value = 1777.2941
$ws.Range("M34").Value = -2310.8813
$ws.Range("N34").Value = -2181.2941

# Row 107
$ws.Range("H107").Value = 557.21875
$ws.Range("I107").Value = 461.2353
$ws.Range("J107").Value = 666
$ws.Range("K107").Value = 461.2353
$ws.Range("L107").Value = 666
$ws.Range("M107").Value = 1458.7647
$ws.Range("N107").Value = -4506

# Row 113
$ws.Range("H113").Value = 1559.2142
$ws.Range("I113").Value = 1102.6364
$ws.Range("J113").Value = 3233.3333
$ws.Range("K113").Value = 1102.6364
$ws.Range("L113").Value = 3233.3333
$ws.Range("M113").Value = 1067.3636
$ws.Range("N113").Value = -7573.3333


# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 767.8099999999999
$ws.Range("I131").Value = 488.66666
$ws.Range("J131").Value = 795.4176
$ws.Range("K131").Value = 1465.99998
$ws.Range("L131").Value = 2386.2528
$ws.Range("M131").Value = 3574.00002
$ws.Range("N131").Value = -12466.2528

# Row 133
$ws.Range("H133").Value = 47623250
$ws.Range("I133").Value = 111112910
$ws.Range("J133").Value = 6000
$ws.Range("K133").Value = 333338730
$ws.Range("L133").Value = 18000
$ws.Range("M133").Value = -333333670
$ws.Range("N133").Value = -28120


# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 133
$ws.Range("H133").Value = 41746.668
$ws.Range("J133").Value = 41746.668
$ws.Range("L133").Value = 41746.668
$ws.Range("N133").Value = -51866.668


# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 52631696
$ws.Range("I55").Value = 118
$ws.Range("J55").Value = 71428690
$ws.Range("K55").Value = 118
$ws.Range("L55").Value = 71428690
$ws.Range("M55").Value = 55
$ws.Range("N55").Value = -71429036

